# Update with scan code and analytical geometry
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column C, matching the header style already used in A1/B1
$ws.Range("C1").Value = "Coord: normal vector scan"
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)  # xlPasteFormats

# Updated angle values (column B) and new normal-vector coordinate strings (column C)
$rows = @(
    @{ Row = 2;  B = 0.4202898136369796;  C = "[0.         0.31802023 0.94808393]" },
    @{ Row = 3;  B = 1.746427337404967;   C = "[-0.43617147  0.51976349  0.73457496]" },
    @{ Row = 4;  B = 0.1630676989240208;  C = "[-9.06280432e-04  1.33439476e-02  9.99910555e-01]" },
    @{ Row = 5;  B = 0.6039428534655938;  C = "[-6.33765013e-04  2.84389540e-01 -9.58708604e-01]" },
    @{ Row = 6;  B = 1.771777786130546;   C = "[0.722527   0.29042097 0.62738376]" },
    @{ Row = 7;  B = 0.8110555212955693;  C = "[-0.73290713 -0.27243498  0.623399  ]" },
    @{ Row = 8;  B = 0.559923347602286;   C = "[0.         0.31570875 0.94885615]" },
    @{ Row = 9;  B = 1.003497847685942;   C = "[ 0.         -0.30835348  0.95127185]" },
    @{ Row = 10; B = 1.887405605205275;   C = "[-0.7295924   0.25409659  0.63492508]" },
    @{ Row = 11; B = 1.877216681045321;   C = "[ 0.72993042 -0.27409819  0.62615634]" },
    @{ Row = 12; B = 0.6323109968505883;  C = "[ 6.54468149e-04 -2.83915167e-01 -9.58849180e-01]" },
    @{ Row = 13; B = 2.486404314623541;   C = "[-0.7161322  -0.28489093  0.63717488]" },
    @{ Row = 14; B = 3.884459959925173;   C = "[0.70064789 0.26004277 0.66443231]" },
    @{ Row = 15; B = 1.458627817372756;   C = "[ 0.         -0.30078739  0.95369122]" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
}
